# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
#
# The source data rows got re-ordered/re-matched (same match id/date groups,
# different row-to-match assignment). For each affected block of rows, every
# column except "A" (the sequential row counter) is rotated one position
# within the block: row[i] (columns B:AD) takes on what used to be
# row[i+1]'s values, with the last row in the block wrapping around to take
# the first row's original values.
#
# Read everything first (so later writes don't clobber values still needed),
# then write the rotated results back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is a contiguous group of worksheet rows that got
# re-matched together. Pairs rotate (= swap); the 291-294 block is a 4-cycle.
$groups = @(
    @(31, 32),
    @(258, 259),
    @(284, 285),
    @(287, 288),
    @(291, 292, 293, 294)
)

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot the B:AD values for every row in this group before writing
    # anything back.
    $snapshots = @()
    foreach ($r in $group) {
        $rng = $ws.Range("B$r`:AD$r")
        $snapshots += ,$rng.Value()
    }

    # row[i] receives what row[i+1] (wrapping) used to hold.
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcValues = $snapshots[($i + 1) % $n]
        $ws.Range("B$destRow`:AD$destRow").Value = $srcValues
    }
}
